$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# URL: http://ibm.com/... -> http://linuxforhealth.org/...
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/patient-importance"

# Version: 7.0.0 -> 8.0.0
$wsMeta.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# Description: "Health Data Connect" -> "LinuxForHealth"
$wsMeta.Range("B11").Value = "Patient importance status codes derived from customer-specific code mappings, used to trigger or limit LinuxForHealth patient operations."

# --- Include from Patient Importan sheet ---
$wsCodes = $wb.Worksheets.Item("Include from Patient Importan")

# System URI: http://ibm.com/... -> http://linuxforhealth.org/...
$wsCodes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/patient-importance"
